$d = $word.ActiveDocument

# 1. Header title "KABUPATEN HULU SUNGAI TENGAH" -> "${u_kabupaten}"
#    (two occurrences, both get the identical replacement)
$d.Content.Find.Execute("KABUPATEN HULU SUNGAI TENGAH", $true, $false, $false, $false, $false, `
    $true, 1, $false, '${u_kabupaten}', 2) | Out-Null

# 2. Signature line " BPS Kabupaten Hulu Sungai Tengah" (preceded by "Kepala") -> " BPS ${kabupaten}"
#    Leading space keeps this distinct from the other "BPS Kabupaten Hulu Sungai Tengah" occurrence
#    and keeps "Kepala" (and its own formatting) untouched.
$d.Content.Find.Execute(" BPS Kabupaten Hulu Sungai Tengah", $true, $false, $false, $false, $false, `
    $true, 1, $false, ' BPS ${kabupaten}', 2) | Out-Null

# 3. Remaining "BPS Kabupaten Hulu Sungai Tengah" (list item, starts the paragraph) -> "BPS ${kabupaten}"
$d.Content.Find.Execute("BPS Kabupaten Hulu Sungai Tengah", $true, $false, $false, $false, $false, `
    $true, 1, $false, 'BPS ${kabupaten}', 2) | Out-Null

# 4. "Barabai" -> "${ibukota}" (two occurrences, both get the identical replacement)
$d.Content.Find.Execute("Barabai", $true, $false, $false, $false, $false, `
    $true, 1, $false, '${ibukota}', 2) | Out-Null
